# Apply cryptos.xlsx price/volume updates (generated from the canonical OOXML diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric need to be forced back to text,
# matching the source file where these columns are plain inline strings
# (e.g. "59.514.39", "5.50", "1.00") rather than numbers. We briefly set a
# text number-format so Excel stores the literal string, then restore the
# cell's original style so no stray formatting is left behind.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '59.514.39'
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').Value = '2.301.71'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue $ws.Range('D5') '539.13'
$ws.Range('E5').Value = '  -1.97%  '
Set-TextValue $ws.Range('D6') '127.61'
$ws.Range('E6').Value = '  -5.34%  '
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue $ws.Range('D8') '0.566'
$ws.Range('E8').Value = '  -4.14%  '
$ws.Range('D9').Value = '2.300.87'
$ws.Range('E9').Value = '  -3.31%  '
$ws.Range('E10').Value = '  -1.93%  '
Set-TextValue $ws.Range('D11') '5.50'
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E13').Value = '  -3.63%  '
$ws.Range('D14').Value = '2.712.09'
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D15') '23.00'
$ws.Range('E15').Value = '  -5.44%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '59.405.55'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('E17').Value = '  -3.15%  '
$ws.Range('D18').Value = '2.329.19'
$ws.Range('E18').Value = '  -2.30%  '
Set-TextValue $ws.Range('D19') '10.36'
$ws.Range('E19').Value = '  -4.56%  '
$ws.Range('E20').Value = '  -5.84%  '
Set-TextValue $ws.Range('D21') '308.53'
$ws.Range('E21').Value = '  -3.47%  '
$ws.Range('E22').Value = '  -6.77%  '
Set-TextValue $ws.Range('D23') '1.00'
$ws.Range('E23').Value = '  -0.46%  '
Set-TextValue $ws.Range('D24') '62.92'
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('E25').Value = '  -3.52%  '
$ws.Range('E26').Value = '  +0.11%  '
Set-TextValue $ws.Range('D27') '7.66'
$ws.Range('E27').Value = '  -6.93%  '
$ws.Range('E28').Value = '  -2.86%  '
Set-TextValue $ws.Range('D29') '171.58'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +2.37%  '
Set-TextValue $ws.Range('D31') '1.69'
$ws.Range('E31').Value = '  -3.65%  '
$ws.Range('D32').Value = '0.0₃0710'
$ws.Range('E32').Value = '  -6.04%  '
Set-TextValue $ws.Range('D33') '5.76'
$ws.Range('E33').Value = '  -3.81%  '
Set-TextValue $ws.Range('D34') '0.376'
$ws.Range('E34').Value = '  -3.28%  '
$ws.Range('E35').Value = '  +0.00%  '
Set-TextValue $ws.Range('D36') '1.31'
$ws.Range('E36').Value = '  -7.70%  '
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  -6.64%  '
Set-TextValue $ws.Range('D40') '309.03'
$ws.Range('E40').Value = '  -5.94%  '
Set-TextValue $ws.Range('D41') '37.59'
$ws.Range('E41').Value = '  -2.49%  '
$ws.Range('E42').Value = '  -6.10%  '
Set-TextValue $ws.Range('D43') '135.27'
$ws.Range('E43').Value = '  -7.86%  '
$ws.Range('E44').Value = '  -3.62%  '
Set-TextValue $ws.Range('D45') '0.0933'
$ws.Range('E45').Value = '  -2.75%  '
Set-TextValue $ws.Range('D46') '0.566'
$ws.Range('E46').Value = '  -0.70%  '
Set-TextValue $ws.Range('D47') '0.0487'
$ws.Range('E47').Value = '  -3.76%  '
Set-TextValue $ws.Range('D48') '18.36'
$ws.Range('E48').Value = '  -6.56%  '
$ws.Range('E49').Value = '  +22.87%  '
Set-TextValue $ws.Range('D50') '0.0210'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('E51').Value = '  -0.53%  '
